# Dividend Calculation.xlsx - update June 2017 dividend figures and
# refresh the saved cell-selection / view state to match.

$wb = $excel.ActiveWorkbook

$yearly = $wb.Worksheets.Item("Yearly")
$allTime = $wb.Worksheets.Item("All Time")

# --- Update the 2017 "June" row (row 8, columns L:N) on the Yearly sheet.
# Column O (Grand Total) and the row-15 / "All Time" sheet totals are all
# driven by SUM formulas, so they recalculate automatically.
$yearly.Range("L8").Value = 118.95
$yearly.Range("M8").Value = 44.98
$yearly.Range("N8").Value = 42.31

# --- Selection on the Yearly sheet moves from O22 to K22. Yearly is
# already the active/tabSelected sheet, so a direct Select() is enough.
$yearly.Range("K22").Select()

# --- Selection (and intended scroll position) on the "All Time" sheet
# moves to A55 / A25. Selecting on a non-active sheet requires activating
# it first; we restore the Yearly sheet as the active tab afterwards so
# tabSelected stays on "Yearly", matching the original workbook.
$allTime.Select()
$allTime.Range("A55").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1

$yearly.Select()
